$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value map for the price/volume refresh (GitHub Actions bot update).
$updates = [ordered]@{
    'D2' = '97.876.56'
    'E2' = '  -1.12%  '
    'D3' = '3.425.32'
    'E3' = '  +3.88%  '
    'E4' = '  -0.05%  '
    'D5' = '255.33'
    'E5' = '  +0.09%  '
    'D6' = '656.21'
    'E6' = '  +4.55%  '
    'E7' = '  +2.32%  '
    'D8' = '0.432'
    'E8' = '  +7.06%  '
    'E9' = '  +8.63%  '
    'E10' = '  +0.01%  '
    'D11' = '3.423.39'
    'E11' = '  +3.92%  '
    'E12' = '  +4.12%  '
    'D13' = '41.97'
    'E13' = '  +2.75%  '
    'E14' = '  +15.22%  '
    'D15' = '0.0000259'
    'E15' = '  +3.34%  '
    'D16' = '97.523.65'
    'D17' = '4.068.42'
    'E17' = '  +3.98%  '
    'D18' = '8.67'
    'E18' = '  +35.05%  '
    'D19' = '3.426.16'
    'E19' = '  +3.73%  '
    'D20' = '17.62'
    'E20' = '  +12.51%  '
    'D21' = '0.506'
    'E21' = '  +48.72%  '
    'D22' = '10.78'
    'E22' = '  +14.23%  '
    'D23' = '3.46'
    'E23' = '  -0.15%  '
    'D24' = '505.88'
    'E24' = '  +3.50%  '
    'E25' = '  +1.44%  '
    'B26' = 'Litecoin'
    'C26' = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
    'D26' = '99.16'
    'E26' = '  +9.61%  '
    'B27' = 'NEARProtocol'
    'C27' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'D27' = '6.13'
    'E27' = '  +6.96%  '
    'D28' = '12.76'
    'E28' = '  +4.68%  '
    'D29' = '3.605.94'
    'E29' = '  +4.10%  '
    'E30' = '  +3.27%  '
    'E31' = '  +6.39%  '
    'D32' = '11.39'
    'E32' = '  +6.17%  '
    'E33' = '  +0.26%  '
    'D34' = '0.998'
    'E34' = '  -0.11%  '
    'D35' = '0.576'
    'E35' = '  +19.17%  '
    'D36' = '29.86'
    'E36' = '  +6.55%  '
    'D37' = '2.28'
    'E37' = '  +16.19%  '
    'D38' = '7.80'
    'E38' = '  +6.27%  '
    'D39' = '1.44'
    'E39' = '  +15.31%  '
    'D40' = '0.155'
    'E40' = '  +2.42%  '
    'D41' = '522.04'
    'E41' = '  +5.15%  '
    'E42' = '  +0.08%  '
    'D43' = '0.875'
    'E43' = '  +11.62%  '
    'D44' = '3.74'
    'E44' = '  -2.31%  '
    'E45' = '  +23.77%  '
    'D46' = '5.60'
    'E46' = '  +15.55%  '
    'D47' = '3.32'
    'E47' = '  +5.06%  '
    'D48' = '8.28'
    'E48' = '  +12.61%  '
    'E49' = '  +0.05%  '
    'D50' = '1.59'
    'E50' = '  +13.99%  '
    'D51' = '2.06'
    'E51' = '  +5.19%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so numeric-looking strings (e.g. "97.876.56") are not
    # reinterpreted as numbers / floats by Excel's input parser, then drop the
    # explicit number-format style again so the cell format matches the original.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.ClearFormats()
}
